# Daily Scrum sheet update: shift weekly entries forward, append new week,
# clear stale remark column, and update view/selection + row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Week header dates (column A) shift forward ---
$ws.Range("A7").Value  = 44259
$ws.Range("A13").Value = 44260
$ws.Range("A20").Value = 44263
$ws.Range("A26").Value = 44264
$ws.Range("A33").Value = 44265
$ws.Range("A39").Value = 44266
$ws.Range("A45").Value = 44267

# --- Week 44259 (rows 8-11): "what was done" / "what will we do" / clear "difficulties" ---
$ws.Range("B8").Value  = "Leitura do enunciado e Análise dos novos Ucs"
$ws.Range("C8").Value  = " Complementação da camada DTO"
$ws.Range("D8").Value  = ""

$ws.Range("B9").Value  = "Implementação UC11/12"
$ws.Range("C9").Value  = "Implementação UC11/12"
$ws.Range("D9").Value  = ""

$ws.Range("B10").Value = "Implementação UC11/12"
$ws.Range("C10").Value = "Implementação UC11/12"
$ws.Range("D10").Value = ""

$ws.Range("B11").Value = "Implementação UC11/12"
$ws.Range("C11").Value = "Implementação UC11/12"
$ws.Range("D11").Value = ""

# --- Week 44260 (rows 14-17) ---
$ws.Range("B14").Value = " Complementação da camada DTO"
$ws.Range("C14").Value = "Desenvolvimento classes de domínio web service"
$ws.Range("D14").Value = ""

$ws.Range("B15").Value = "Implementação UC11/12"
$ws.Range("C15").Value = "Integração com BD"
$ws.Range("D15").Value = ""

$ws.Range("B16").Value = "Implementação UC11/12"
$ws.Range("C16").Value = "Integração com BD"
$ws.Range("D16").Value = ""

$ws.Range("B17").Value = "Implementação UC11/12"
$ws.Range("C17").Value = "Desenvolvimento classes de domínio web service"
$ws.Range("D17").Value = ""

# --- Week 44263 (rows 21-24) ---
$ws.Range("B21").Value = "Desenvolvimento classes de domínio web service"
$ws.Range("C21").Value = "Camada DTO/Domain"
$ws.Range("D21").Value = ""

$ws.Range("B22").Value = "Integração com BD"
$ws.Range("C22").Value = "Camada Repositório"
$ws.Range("D22").Value = ""

$ws.Range("B23").Value = "Integração com BD"
$ws.Range("C23").Value = "Camada Controller/Service"
$ws.Range("D23").Value = ""

$ws.Range("B24").Value = "Desenvolvimento classes de domínio web service"
$ws.Range("C24").Value = "Camada Service/Domain"
$ws.Range("D24").Value = ""

# Rows 21-24 got taller to fit new text
$ws.Rows.Item(21).RowHeight = 83.25
$ws.Rows.Item(22).RowHeight = 83.25
$ws.Rows.Item(23).RowHeight = 83.25
$ws.Rows.Item(24).RowHeight = 83.25

# --- Week 44264 (rows 27-30) ---
$ws.Range("B27").Value = "Camada DTO/Domain"
$ws.Range("C27").Value = "Finalizações e testes do webservice"
$ws.Range("D27").Value = ""

$ws.Range("B28").Value = "Camada Repositório"
$ws.Range("C28").Value = "Finalizações e testes do webservice"
$ws.Range("D28").Value = ""

$ws.Range("B29").Value = "Camada Controller/Service"
$ws.Range("C29").Value = "Finalizações e testes do webservice"
$ws.Range("D29").Value = ""

$ws.Range("B30").Value = "Camada Service/Domain"
$ws.Range("C30").Value = "Finalizações e testes do webservice"
$ws.Range("D30").Value = ""

# --- Week 44265 (rows 34-37) ---
$ws.Range("B34").Value = "Finalizações e testes do webservice"
$ws.Range("C34").Value = "Javadocs, Alterações no projecto para integrar novo web-service"

$ws.Range("B35").Value = "Finalizações e testes do webservice"
$ws.Range("C35").Value = "Alterações no projecto para integrar novo web-service"

$ws.Range("B36").Value = "Finalizações e testes do webservice"
$ws.Range("C36").Value = "Alterações no projecto para integrar novo web-service"

$ws.Range("B37").Value = "Finalizações e testes do webservice"
$ws.Range("C37").Value = "Javadocs, Alterações no projecto para integrar novo web-service"

# --- Week 44266 (rows 40-43) ---
$ws.Range("B40").Value = "Javadocs, Alterações no projecto para integrar novo web-service"
$ws.Range("C40").Value = "Testes funcionais à aplicação, correção de erros"

$ws.Range("B41").Value = "Alterações no projecto para integrar novo web-service"
$ws.Range("C41").Value = "Testes funcionais à aplicação, correção de erros"

$ws.Range("B42").Value = "Alterações no projecto para integrar novo web-service"
$ws.Range("C42").Value = "Testes funcionais à aplicação, correção de erros"

$ws.Range("B43").Value = "Javadocs, Alterações no projecto para integrar novo web-service"
$ws.Range("C43").Value = "Testes funcionais à aplicação, correção de erros"

# --- Week 44267 (rows 46-49) ---
$ws.Range("B46").Value = "Testes funcionais à aplicação, correção de erros"
$ws.Range("C46").Value = "Finalização do trabalho para entrega"

$ws.Range("B47").Value = "Testes funcionais à aplicação, correção de erros"
$ws.Range("C47").Value = "Finalização do trabalho para entrega"

$ws.Range("B48").Value = "Testes funcionais à aplicação, correção de erros"
$ws.Range("C48").Value = "Finalização do trabalho para entrega"

$ws.Range("B49").Value = "Testes funcionais à aplicação, correção de erros"
$ws.Range("C49").Value = "Finalização do trabalho para entrega"

# --- View state: scroll down to the new last week, select D47 ---
$ws.Range("D47").Select()
$excel.ActiveWindow.ScrollRow = 43
